# Update the yearly cost sheet: drop the oldest reporting year
# (1396/12) and append the newest one (1401/12). Every 5-year data
# column (E:I) shifts one column to the left (E<-F, F<-G, G<-H,
# H<-I) and the freed-up I column receives the new year's figure.
# This mirrors the "read_price" recompute that rolled the trailing
# twelve-month window forward by one fiscal year.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Year-label header rows: shift the "Twelve months ended ..." ---
#        captions left by one column, then stamp the new year label
#        into column I.
$headerRows = @(8, 27, 35, 43, 51, 59, 67, 75, 83, 91, 98, 105, 112, 119)
foreach ($r in $headerRows) {
    $ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
    $ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
    $ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 8).Value = $ws.Cells.Item($r, 9).Value2
    $ws.Cells.Item($r, 9).Value = "دوازده ماهه منتهی به 1401/12"
}

# --- 2. Data rows: shift the five trailing-year figures left and ---
#        drop in the freshly computed value for the new year (I).
$newYearValue = @{
    10 = 2732352649
    11 = 2542709
    12 = 42759866
    13 = 2777655224
    14 = 0
    15 = 2777655224
    16 = 4436333
    17 = 0
    18 = 2782091557
    19 = 141839344
    20 = -161386250
    21 = 2762544651
    22 = 0
    23 = 2762544651
    29 = 244189
    30 = 3628446
    31 = 3872635
    37 = 8758853
    38 = 116907874
    39 = 125666727
    45 = 8782050
    46 = 117000277
    47 = 125782327
    53 = 220992
    54 = 3536043
    55 = 3757035
    61 = 3624293
    62 = 68520592
    63 = 72144885
    69 = 190429286
    70 = 2548671338
    71 = 2739100624
    77 = 189290268
    78 = 2543062381
    79 = 2732352649
    85 = 4763311
    86 = 74129549
    87 = 78892860
    93 = 14842163
    94 = 18884280
    100 = 21741350
    101 = 21800682
    107 = 21554223
    108 = 21735524
    114 = 21554224
    115 = 20963984
    121 = 0
    122 = 0
    123 = 0
    124 = 0
    125 = 1966042
    126 = 0
    127 = 3987989
    128 = 11933627
    129 = 0
    130 = 24872208
    131 = 42759866
}

foreach ($r in $newYearValue.Keys) {
    $f = $ws.Cells.Item($r, 6).Value2
    $g = $ws.Cells.Item($r, 7).Value2
    $h = $ws.Cells.Item($r, 8).Value2
    $i = $ws.Cells.Item($r, 9).Value2

    $ws.Cells.Item($r, 5).Value = $f
    $ws.Cells.Item($r, 6).Value = $g
    $ws.Cells.Item($r, 7).Value = $h
    $ws.Cells.Item($r, 8).Value = $i
    $ws.Cells.Item($r, 9).Value = $newYearValue[$r]
}
